$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# Update the keyword text in A10 (was "verify signup link")
$ws.Range("A10").Value = "verify forgot password link"

# Update the view: zoom to 175% and move the selection to B16
$ws.Activate()
$excel.ActiveWindow.Zoom = 175
$ws.Range("B16").Select()
